# Update "想去人数" (want-to-go count) figures that were refreshed in the
# gh-pages data regeneration (commit 456a3b4).
#
# Affected sheets/cells:
#   展览 (Exhibition):   F9 203->204, F20 5358->5359, F23 787->790, F25 265->266
#   全部类型 (All types): F9 203->204, F21 5358->5359, F25 787->790, F27 265->266

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F9").Value = 204
$wsExhibition.Range("F20").Value = 5359
$wsExhibition.Range("F23").Value = 790
$wsExhibition.Range("F25").Value = 266

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 204
$wsAll.Range("F21").Value = 5359
$wsAll.Range("F25").Value = 790
$wsAll.Range("F27").Value = 266
